$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-34 down to 5-35
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the new data
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44537
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103003
$ws.Range("J4").Value = "Damasco"
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S4").Value = 889
$ws.Range("T4").Value = 18
